# Add 2022-Q4 data: insert a new quarter sheet (cloned from 2022-Q3) and
# update the "总计" (summary) sheet with the new quarter's row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by cloning the existing "2022-Q3"
#    sheet (same columns/headers/formatting), inserted right before it.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3, $null)
$wsQ4 = $wb.Worksheets.Item("总计").Next
$wsQ4.Name = "2022-Q4"

# Update the cloned sheet's figures for the Q4 snapshot.
$wsQ4.Range("D2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "14.15"
$wsQ4.Range("E2").NumberFormat = "@"
$wsQ4.Range("E2").Value = "75.21"
$wsQ4.Range("F2").NumberFormat = "@"
$wsQ4.Range("F2").Value = "3.15"
$wsQ4.Range("G2").NumberFormat = "@"
$wsQ4.Range("G2").Value = "0.4457"
$wsQ4.Range("H2").Value = 9

$wsQ4.Range("E3").NumberFormat = "@"
$wsQ4.Range("E3").Value = "75.21"
$wsQ4.Range("F3").NumberFormat = "@"
$wsQ4.Range("F3").Value = "3.15"
$wsQ4.Range("G3").NumberFormat = "@"
$wsQ4.Range("G3").Value = "0.0123"
$wsQ4.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new row for 2022-Q4 at the top of
#    the data (row 2), pushing the existing Q3/Q2 rows down.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Carry the "A" column style down from the row below onto the new row.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.46

# The row-insert shifts the old A3/A4 index values down too; restore the
# sequential index column (0, 1, 2, ...) to match the target data.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
